# Advent of Code 2024, Day 8 - append the new day's runtime result and
# correct the runtimes that changed for days 1-4.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Corrected runtimes for days 1-4 (column B, rows 3-6)
$ws.Range("B3").Value = 0.0023924599999999999
$ws.Range("B4").Value = 0.0034427400000000001
$ws.Range("B5").Value = 0.0022166
$ws.Range("B6").Value = 0.0043159599999999998

# New rows for days 5-8
$ws.Range("A7").Value = 5
$ws.Range("B7").Value = 0.0092547800000000006

$ws.Range("A8").Value = 6
$ws.Range("B8").Value = 16.177439159999999

$ws.Range("A9").Value = 7
$ws.Range("B9").Value = 2.1248556199999999

$ws.Range("A10").Value = 8
$ws.Range("B10").Value = 0.0033554600000000002

# Match the author's on-disk selection (A3:B10, matching the new data extent)
[void]$ws.Range("A3:B10").Select()
